$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

function Set-PlainValue($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# Row 2
Set-TextValue 2 4 '41.488.58'
Set-PlainValue 2 5 '  -2.45%  '

# Row 3
Set-TextValue 3 4 '2.477.00'
Set-PlainValue 3 5 '  -1.82%  '

# Row 4
Set-PlainValue 4 5 '  +0.78%  '

# Row 5
Set-TextValue 5 4 '313.25'
Set-PlainValue 5 5 '  -0.41%  '

# Row 6
Set-TextValue 6 4 '92.63'
Set-PlainValue 6 5 '  -6.34%  '

# Row 7
Set-TextValue 7 4 '0.545'
Set-PlainValue 7 5 '  -3.16%  '

# Row 8
Set-PlainValue 8 5 '  +0.68%  '

# Row 9
Set-PlainValue 9 5 '  -4.46%  '

# Row 10
Set-TextValue 10 4 '33.09'
Set-PlainValue 10 5 '  -5.96%  '

# Row 11
Set-TextValue 11 4 '0.0780'
Set-PlainValue 11 5 '  -2.63%  '

# Row 12
Set-PlainValue 12 5 '  -0.14%  '

# Row 13
Set-TextValue 13 4 '2.861.49'
Set-PlainValue 13 5 '  -1.70%  '

# Row 14
Set-PlainValue 14 5 '  -4.65%  '

# Row 15
Set-PlainValue 15 2 'WrappedEther'
Set-PlainValue 15 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 15 4 '2.484.65'
Set-PlainValue 15 5 '  +0.11%  '

# Row 16
Set-PlainValue 16 2 'Chainlink'
Set-PlainValue 16 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 16 4 '15.33'
Set-PlainValue 16 5 '  +0.80%  '

# Row 17
Set-TextValue 17 4 '0.784'
Set-PlainValue 17 5 '  -3.08%  '

# Row 18
Set-TextValue 18 4 '41.340.14'
Set-PlainValue 18 5 '  -2.79%  '

# Row 19
Set-PlainValue 19 2 'Uniswap'
Set-PlainValue 19 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 19 4 '6.28'
Set-PlainValue 19 5 '  -4.72%  '

# Row 20
Set-PlainValue 20 2 'ShibaInu'
Set-PlainValue 20 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 20 4 '0.0₃0923'
Set-PlainValue 20 5 '  -1.72%  '

# Row 21
Set-TextValue 21 4 '69.88'
Set-PlainValue 21 5 '  +1.23%  '

# Row 22
Set-TextValue 22 4 '11.03'
Set-PlainValue 22 5 '  -9.11%  '

# Row 23
Set-TextValue 23 4 '234.66'
Set-PlainValue 23 5 '  -2.83%  '

# Row 24
Set-PlainValue 24 5 '  -4.08%  '

# Row 25
Set-PlainValue 25 5 '  -0.06%  '

# Row 26
Set-TextValue 26 4 '1.87'
Set-PlainValue 26 5 '  -5.71%  '

# Row 27
Set-TextValue 27 4 '23.99'
Set-PlainValue 27 5 '  -5.96%  '

# Row 28
Set-PlainValue 28 5 '  -0.27%  '

# Row 29
Set-TextValue 29 4 '9.73'
Set-PlainValue 29 5 '  -2.58%  '

# Row 30
Set-TextValue 30 4 '36.44'
Set-PlainValue 30 5 '  -3.31%  '

# Row 31
Set-TextValue 31 4 '152.95'
Set-PlainValue 31 5 '  -2.10%  '

# Row 32
Set-TextValue 32 4 '5.44'
Set-PlainValue 32 5 '  -8.08%  '

# Row 33
Set-PlainValue 33 5 '  -3.00%  '

# Row 34
Set-TextValue 34 4 '2.53'
Set-PlainValue 34 5 '  -6.43%  '

# Row 35
Set-TextValue 35 4 '0.0747'
Set-PlainValue 35 5 '  -4.43%  '

# Row 36
Set-TextValue 36 4 '17.72'
Set-PlainValue 36 5 '  +1.21%  '

# Row 37
Set-PlainValue 37 5 '  -3.86%  '

# Row 38
Set-PlainValue 38 5 '  -5.82%  '

# Row 39
Set-PlainValue 39 5 '  -3.38%  '

# Row 40
Set-PlainValue 40 5 '  -7.28%  '

# Row 41
Set-TextValue 41 4 '4.03'
Set-PlainValue 41 5 '  -4.45%  '

# Row 42
Set-PlainValue 42 5 '  +1.03%  '

# Row 43
Set-TextValue 43 4 '19.63'
Set-PlainValue 43 5 '  -11.22%  '

# Row 44
Set-TextValue 44 4 '1.964.10'
Set-PlainValue 44 5 '  -2.60%  '

# Row 45
Set-PlainValue 45 5 '  -4.53%  '

# Row 46
Set-TextValue 46 4 '2.95'
Set-PlainValue 46 5 '  -8.20%  '

# Row 47
Set-TextValue 47 4 '8.76'
Set-PlainValue 47 5 '  -1.98%  '

# Row 48
Set-TextValue 48 4 '2.724.13'
Set-PlainValue 48 5 '  -1.50%  '

# Row 49
Set-PlainValue 49 2 'ordi'
Set-PlainValue 49 3 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue 49 4 '68.28'
Set-PlainValue 49 5 '  -4.42%  '

# Row 50
Set-PlainValue 50 2 'Aave'
Set-PlainValue 50 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 50 4 '95.89'
Set-PlainValue 50 5 '  -4.08%  '

# Row 51
Set-TextValue 51 4 '0.175'
Set-PlainValue 51 5 '  -6.47%  '

